$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-29 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-30 Saturday", 2) | Out-Null
$d.Content.Find.Execute("20×57=1140", $true, $false, $false, $false, $false, $true, 1, $false, "38×81=3078", 2) | Out-Null
$d.Content.Find.Execute("84×81=6804", $true, $false, $false, $false, $false, $true, 1, $false, "60×75=4500", 2) | Out-Null
$d.Content.Find.Execute("79×50=3950", $true, $false, $false, $false, $false, $true, 1, $false, "12×34=408", 2) | Out-Null
$d.Content.Find.Execute("39×76=2964", $true, $false, $false, $false, $false, $true, 1, $false, "90×26=2340", 2) | Out-Null
$d.Content.Find.Execute("68×25=1700", $true, $false, $false, $false, $false, $true, 1, $false, "60×44=2640", 2) | Out-Null
$d.Content.Find.Execute("18×22=396", $true, $false, $false, $false, $false, $true, 1, $false, "92×73=6716", 2) | Out-Null
$d.Content.Find.Execute("21×68=1428", $true, $false, $false, $false, $false, $true, 1, $false, "52×13=676", 2) | Out-Null
$d.Content.Find.Execute("70×61=4270", $true, $false, $false, $false, $false, $true, 1, $false, "30×55=1650", 2) | Out-Null
$d.Content.Find.Execute("89×75=6675", $true, $false, $false, $false, $false, $true, 1, $false, "74×52=3848", 2) | Out-Null
$d.Content.Find.Execute("87×32=2784", $true, $false, $false, $false, $false, $true, 1, $false, "53×95=5035", 2) | Out-Null
$d.Content.Find.Execute("68×99=6732", $true, $false, $false, $false, $false, $true, 1, $false, "51×52=2652", 2) | Out-Null
$d.Content.Find.Execute("40×80=3200", $true, $false, $false, $false, $false, $true, 1, $false, "44×32=1408", 2) | Out-Null
$d.Content.Find.Execute("60×74=4440", $true, $false, $false, $false, $false, $true, 1, $false, "57×11=627", 2) | Out-Null
$d.Content.Find.Execute("40×39=1560", $true, $false, $false, $false, $false, $true, 1, $false, "99×40=3960", 2) | Out-Null
$d.Content.Find.Execute("40×72=2880", $true, $false, $false, $false, $false, $true, 1, $false, "46×54=2484", 2) | Out-Null
$d.Content.Find.Execute("54×33=1782", $true, $false, $false, $false, $false, $true, 1, $false, "66×33=2178", 2) | Out-Null
$d.Content.Find.Execute("95×63=5985", $true, $false, $false, $false, $false, $true, 1, $false, "26×75=1950", 2) | Out-Null
$d.Content.Find.Execute("83×68=5644", $true, $false, $false, $false, $false, $true, 1, $false, "38×72=2736", 2) | Out-Null
$d.Content.Find.Execute("44×58=2552", $true, $false, $false, $false, $false, $true, 1, $false, "77×25=1925", 2) | Out-Null
$d.Content.Find.Execute("23×63=1449", $true, $false, $false, $false, $false, $true, 1, $false, "83×63=5229", 2) | Out-Null
$d.Content.Find.Execute("62×43=2666", $true, $false, $false, $false, $false, $true, 1, $false, "78×14=1092", 2) | Out-Null
$d.Content.Find.Execute("47×17=799", $true, $false, $false, $false, $false, $true, 1, $false, "40×85=3400", 2) | Out-Null
$d.Content.Find.Execute("79×69=5451", $true, $false, $false, $false, $false, $true, 1, $false, "84×18=1512", 2) | Out-Null
$d.Content.Find.Execute("72×88=6336", $true, $false, $false, $false, $false, $true, 1, $false, "75×41=3075", 2) | Out-Null
$d.Content.Find.Execute("98×61=5978", $true, $false, $false, $false, $false, $true, 1, $false, "41×77=3157", 2) | Out-Null
